# RPA datasets push 2024-05-14
# Update the "38커뮤니케이션(최근일자기준)" IPO tracker sheet:
#  - add two newly announced deals (KB스팩29호, 에이치엠씨아이비스팩7호)
#    right after the first (most recent) row
#  - drop the two oldest rows that rolled off the bottom of the window
#    (코칩, 유안타스팩16호)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# Insert two blank rows right after row 2 (shifts the rest of the table down)
$ws.Range("3:4").Insert()

# Row 3: KB스팩29호
$ws.Cells.Item(3, 1).Value = "KB스팩29호"
$ws.Cells.Item(3, 2).Value = "2024.06.04~06.05"
$ws.Cells.Item(3, 3).Value = "2,000~2,000"
$ws.Cells.Item(3, 4).Value = "-"
$ws.Cells.Item(3, 5).Value = "12000"
$ws.Cells.Item(3, 6).Value = "KB증권"

# Row 4: 에이치엠씨아이비스팩7호
$ws.Cells.Item(4, 1).Value = "에이치엠씨아이비스팩7호"
$ws.Cells.Item(4, 2).Value = "2024.06.04~06.05"
$ws.Cells.Item(4, 3).Value = "2,000~2,000"
$ws.Cells.Item(4, 4).Value = "-"
$ws.Cells.Item(4, 5).Value = "14000"
$ws.Cells.Item(4, 6).Value = "현대차증권"

# The two inserted rows pushed the table to 22 rows (1 header + 21 data);
# drop the two oldest data rows now sitting at the bottom (22:23)
$ws.Range("22:23").Delete()
